$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'306.19"
$ws.Range("E2").Value = "'0.09%"
$ws.Range("D3").Value = "'36.37"
$ws.Range("E3").Value = "'-0.78%"
$ws.Range("D4").Value = "'5.046"
$ws.Range("E4").Value = "'0.33%"
$ws.Range("D5").Value = "'0.07899"
$ws.Range("E5").Value = "'0.69%"
$ws.Range("D6").Value = "'2.131"
$ws.Range("E6").Value = "'-1.51%"
$ws.Range("B7").Value = "GateToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D7").Value = "'4.160"
$ws.Range("E7").Value = "'2.45%"
$ws.Range("B8").Value = "KuCoinToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
$ws.Range("D8").Value = "'7.976"
$ws.Range("E8").Value = "'-0.77%"
$ws.Range("B9").Value = "MXToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D9").Value = "'0.9269"
$ws.Range("E9").Value = "'0.53%"
$ws.Range("B10").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C10").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D10").Value = "'0.09706"
$ws.Range("E10").Value = "'-2.77%"
$ws.Range("B11").Value = "WazirX"
$ws.Range("C11").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D11").Value = "'0.1863"
$ws.Range("E11").Value = "'-0.68%"
$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D12").Value = "'0.09061"
$ws.Range("E12").Value = "'3.55%"
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D13").Value = "'0.03707"
$ws.Range("E13").Value = "'2.63%"
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D14").Value = "'0.09907"
$ws.Range("E14").Value = "'-0.33%"
$ws.Range("B15").Value = "BitForexToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D15").Value = "'0.001441"
$ws.Range("E15").Value = "'-4.17%"
$ws.Range("B16").Value = "TigerCash"
$ws.Range("C16").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D16").Value = "'0.005671"
$ws.Range("E16").Value = "'0.07%"
$ws.Range("B17").Value = "LEO"
$ws.Range("C17").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D17").Value = "'3.468"
$ws.Range("E17").Value = "'0.26%"
$ws.Range("E18").Value = "'13.26%"
$ws.Range("D19").Value = "'0.3422"
$ws.Range("E19").Value = "'-0.84%"
$ws.Range("D20").Value = "'0.1312"
$ws.Range("E20").Value = "'-2.67%"
$ws.Range("D21").Value = "'5.129"
$ws.Range("E21").Value = "'4.24%"
$ws.Range("D22").Value = "'0.2257"
$ws.Range("E22").Value = "'2.41%"
$ws.Range("D23").Value = "'0.04579"
$ws.Range("E23").Value = "'-0.79%"
$ws.Range("D24").Value = "'0.001238"
$ws.Range("E24").Value = "'0.37%"
$ws.Range("D25").Value = "'0.004794"
$ws.Range("E25").Value = "'-7.65%"
$ws.Range("D26").Value = "'0.0001305"
$ws.Range("E26").Value = "'-6.80%"
$ws.Range("E27").Value = "'74.08%"
$ws.Range("D39").Value = "'0.01957"
$ws.Range("E39").Value = "'8.08%"
$ws.Range("D40").Value = "'0.04926"
$ws.Range("E40").Value = "'3.65%"
$ws.Range("D41").Value = "'0.007743"
$ws.Range("E41").Value = "'-2.10%"
$ws.Range("D42").Value = "'0.1396"
$ws.Range("E42").Value = "'-0.83%"
$ws.Range("D43").Value = "'0.007840"
$ws.Range("E43").Value = "'3.05%"
$ws.Range("D44").Value = "'0.002149"
$ws.Range("E44").Value = "'-1.46%"
$ws.Range("E45").Value = "'11.66%"
$ws.Range("D46").Value = "'0.00006303"
$ws.Range("E46").Value = "'-0.94%"
$ws.Range("D47").Value = "'0.00000000753"
$ws.Range("E47").Value = "'0.23%"
$ws.Range("E48").Value = "'0.04%"
$ws.Range("D49").Value = "'51.71"
$ws.Range("E49").Value = "'43.06%"
$ws.Range("D50").Value = "'0.001908"
$ws.Range("E50").Value = "'-29.16%"
$ws.Range("D51").Value = "'0.00002109"
$ws.Range("E51").Value = "'0.23%"
